# Fix Training Data Issue
# The "Date" column was populated with the wrong string format (e.g. "6-17-2013-14"
# instead of the intended ISO-ish "2014-06-17"), because of how the NBA stats site
# displayed the game date. This rewrites every "Date" cell in the sheet that still
# has the bad value to the corrected one, without letting Excel reinterpret the
# replacement text as a real date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-17-2013-14"
$newDate = "2014-06-17"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Locate the "Date" column header on row 1.
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
$dateHeaderCell = $headerRange.Find("Date")
if ($dateHeaderCell -eq $null) {
    $dateCol = 58
} else {
    $dateCol = $dateHeaderCell.Column
}

$dataRange = $ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol))

# Force text storage so Excel doesn't silently convert "2014-06-17" into a date
# serial number (and reformat the cell) when we write it back.
$dataRange.NumberFormat = "@"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldDate) {
        $cell.Value2 = $newDate
    }
}

# Restore the default style so touched cells keep looking exactly like before
# (only their text content changed).
$dataRange.Style = "Normal"
